# Update the "想去人数" (want-to-go headcount) column F values across sheets
# to reflect regenerated site data (gh-pages output at 456a3b4).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("F9").Value = 7416
$ws.Range("F10").Value = 486
$ws.Range("F13").Value = 35
$ws.Range("F14").Value = 148
$ws.Range("F16").Value = 229
$ws.Range("F18").Value = 1328
$ws.Range("F24").Value = 0
$ws.Range("F26").Value = 0
$ws.Range("F27").Value = 21
$ws.Range("F28").Value = 0
$ws.Range("F29").Value = 166
$ws.Range("F30").Value = 5218
$ws.Range("F31").Value = 0
$ws.Range("F32").Value = 0
$ws.Range("F33").Value = 2761
$ws.Range("F34").Value = 0
$ws.Range("F36").Value = 12
$ws.Range("F37").Value = 0
$ws.Range("F40").Value = 0
$ws.Range("F42").Value = 248
$ws.Range("F44").Value = 3978

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 2

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 0
$ws.Range("F4").Value = 19581
$ws.Range("F7").Value = 1089
$ws.Range("F9").Value = 0
$ws.Range("F10").Value = 0
$ws.Range("F11").Value = 0
$ws.Range("F14").Value = 148
$ws.Range("F15").Value = 0
$ws.Range("F19").Value = 367
$ws.Range("F23").Value = 0
$ws.Range("F24").Value = 60
$ws.Range("F26").Value = 1070
$ws.Range("F27").Value = 21
$ws.Range("F31").Value = 0
$ws.Range("F33").Value = 49
$ws.Range("F34").Value = 33
$ws.Range("F35").Value = 0
$ws.Range("F37").Value = 84
$ws.Range("F39").Value = 0
$ws.Range("F40").Value = 1320
$ws.Range("F41").Value = 0
$ws.Range("F42").Value = 0
$ws.Range("F45").Value = 0
$ws.Range("F46").Value = 3978
